# Apply updates for "Add data for 2022-04-11"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-03"

# Update the header label in I1 (shared string) to match the new date
$ws.Range("I1").Value = "2022 (through 04-03)"

# Update the April total (row 5, since row1=header, row2=Jan,... row5=April)
$ws.Range("I5").Value = 7

# Update the overall Total row (row 14) for the Total column (I)
$ws.Range("I14").Value = 440
